# Change table placeholder syntax from "${table:name.field}" to
# "${table:name:field}" so multi-dimensional key names can be supported.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "B7"  = '${table:ages:name}'
    "C7"  = '${table:ages:age}'
    "E7"  = '${table:scores:name}'
    "F7"  = '${table:scores:score}'
    "B10" = '${table:coords:x}'
    "C10" = '${table:coords:y}'
    "B13" = '${table:dates:name}'
    "C13" = '${table:dates:dates}'
}

foreach ($addr in $replacements.Keys) {
    $ws.Range($addr).Value = $replacements[$addr]
}
